# Apply the edit described by the diff:
#  - Insert a new worksheet "_h2_assigned_courses__h2_ (39)" before "feedback",
#    populated with a small table of assigned-course enrollment data.
#  - Re-activate the "feedback" sheet (it stays the visually-selected tab),
#    restoring its scroll/selection state near the bottom of the sheet.
#  - Touch a couple of the summary-row formulas on "feedback" (AVERAGE/STDEV)
#    that were re-entered by the author (no value change, just re-entry).

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new sheet right before "feedback" ------------------------
$ws = $wb.Worksheets.Add($wb.Worksheets.Item("feedback"))
$ws.Name = "_h2_assigned_courses__h2_ (39)"

# Re-fetch "feedback" after the Add() call -- a reference captured before a
# Worksheets-collection mutation doesn't track the live sheet afterwards.
$fb = $wb.Worksheets.Item("feedback")

# --- 2. Fill in the data ------------------------------------------------
# The order in which new distinct strings are first written controls the
# order they land in in the shared-string table, so write the handful of
# repeated text values first (Subject/Course-Number text values), then the
# header row right-to-left, then backfill the remaining plain numbers.

$ws.Range("B2").Value = "CS"
$ws.Range("C9").Value = "890DW"
$ws.Range("C8").Value = "890DE"
$ws.Range("C5").Value = "490CN"
$ws.Range("C4").Value = "490AB"

$ws.Range("E1").Value = "Enrollment"
$ws.Range("D1").Value = "Section Number"
$ws.Range("C1").Value = "Course Number"
$ws.Range("B1").Value = "Subject Code"
$ws.Range("A1").Value = "CRN"

# CRN (A), remaining Subject (B), remaining Course Number (C),
# Section Number (D) and Enrollment (E) columns for rows 2-15.
$data = @(
    @(10667, "CS", 205, 1, 9),
    @(10680, "CS", 280, 1, 53),
    @(13499, "CS", "490AB", 1, 1),
    @(13495, "CS", "490CN", 1, 1),
    @(10705, "CS", 499, 1, 1),
    @(10706, "CS", 499, 2, 1),
    @(13229, "CS", "890DE", 1, 1),
    @(13496, "CS", "890DW", 1, 2),
    @(10716, "CS", 900, 1, 3),
    @(10717, "CS", 900, 2, 2),
    @(10718, "CS", 900, 3, 1),
    @(10719, "CS", 900, 4, 0),
    @(10723, "CS", 901, 5, 2),
    @(10738, "CS", 902, 5, 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# --- 3. Match page-margins used elsewhere in the workbook -------------------
$ps = $ws.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

# --- 4. Touch a few formulas on "feedback" that were re-entered -------------
$fb.Range("C77").Formula = "=AVERAGE(C2:C76)"
$fb.Range("D77").Formula = "=AVERAGE(D2:D76)"
$fb.Range("C78").Formula = "=STDEV(C2:C6)"

# --- 5. Re-select "feedback" as the active sheet, restoring scroll/selection -
$fb.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1
$fb.Range("D77").Select()
